# Rename the worksheet/sheet tab to the new report title.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "OS OPRs"

# Swap columns B and C: both their header text (B1/C1) and their widths.
$bWidth = $ws.Columns.Item(2).ColumnWidth
$cWidth = $ws.Columns.Item(3).ColumnWidth

$bHeader = $ws.Range("B1").Value2
$cHeader = $ws.Range("C1").Value2

$ws.Range("B1").Value2 = $cHeader
$ws.Range("C1").Value2 = $bHeader

$ws.Columns.Item(2).ColumnWidth = $cWidth
$ws.Columns.Item(3).ColumnWidth = $bWidth

# Update the active cell selection in the frozen-pane view to A2.
$ws.Range("A2").Select()
